$wb = $excel.ActiveWorkbook

# --- Update simulation result lists (append Week 15 data) ---
$ws = $wb.Worksheets.Item('YDS')
$ws.Range('B2').Value = '-3 -1 2 7 -1 2 2 0 0 1 16 1 19 4 8 8 -1 4 7 8 1 2 1 -2 2 7 0 1 4 7 3 2 4 8 4 1 11 3 4 -2 9 2 6 7 0 60 3 6 2 19 0 3 8 4 4 1 0 12 4 1 4 1 3 4 17 0 19 1 -1 4 2 1 3 4 8 2 8 2 9 5 1 2 3 2 28 10 9 2 2 12 1 5 1 -1 13 0 7 4 0 2 2 2 4 16 1 7 2 0 4 3 5 2 3 1 2 11 2 3 18 9 11 22 1 2 2 7 5 5 1 0 -4 2 4 4 4 6 -1 0 1 5 9 5 11 1 6 9 4 2 3 3 8 7 13 9 1 4 -2 0 4 15 0 9 5 4 11 10 -1 2 4 3 76 2 -3 7 4 6 3 19 0 4 3 1 2 3 0 -1 13 3 9 2 3 4 0 8 5 0 3 4 -1 2 4 0 11 1 -1 2 2 4 6 2 1 2 8 5 4 -4 9 3 2 -4 5 3 1 3 4 1 13 2 5 9 2 1 0 6 1 1 2 2 -1 4 2 4 5 3 3 6 -1 -3 2 3 1 3 -3 8 -7 0 1 4 2 14 2 11 6 3 4 2 1 6 0 3 2 1 1 2 4 2 5 4 -1 6 5 7 1 2 1 5 3 2 2 1 1 0 2 3 4 5 6 -5 0 1 11 7 0 16 4 12 5 4 9 -1 8 -1 9 -1 1 2 16 2 2 1 -3 -1 4 5 2 11 2 5 2 1 4 2 17 3 13 3 68 6 30 5 6 3 7 14 4 2 8 4 3 1 2 1 5 12 9 2 3 1 2 6 3 17 0 5 4 4 6 8 2 5 3 4 3 2 -5 3 1 1 -1 2 4 5 6 1 1 12 0 2 -1'

$ws = $wb.Worksheets.Item('YDS')
$ws.Range('B3').Value = '3 4 8 12 11 39 10 9 9 4 4 19 13 9 10 5 10 8 10 8 7 10 6 14 13 5 17 51 15 15 19 6 15 15 22 12 12 5 15 16 10 14 4 9 12 4 5 6 6 6 9 3 25 7 18 26 14 13 13 -2 10 11 13 7 3 27 5 11 14 16 4 7 10 0 12 7 10 9 1 23 9 14 18 9 15 12 11 6 6 2 20 -6 5 9 12 10 12 25 10 22 22 8 14 26 20 6 7 1 14 11 5 48 7 8 8 8 4 13 20 6 6 14 22 5 8 9 10 18 3 5 3 24 8 13 22 10 2 20 6 12 9 11 6 5 8 24 8 7 3 8 8 6 14 7 57 13 10 11 7 5 10 11 14 9 10 10 13 11 13 15 5 -3 7 2 10 14 1 8 3 5 9 9 4 7 11 7 13 16 3 2 50 9 16 17 10 1 16 16 9 10 -1 2 5 7 2 39 0 11 10 10 11 15 1 8 7 12 7 10 8 12 5 11 -1 16 9 5 7 18 -2 4 7 46 11 5 9 9 4 4 8 14 7 8 10 2 14 8 1 2 3 24 7 14 9 6 11 15 1 16 9 12 5 4 8 11 4 10 24 9 5 13 7 12'

$ws = $wb.Worksheets.Item('YDS')
$ws.Range('C2').Value = '10 0 4 3 10 7 11 8 4 8 0 0 6 2 4 3 7 16 4 2 6 1 5 3 2 1 2 6 8 2 1 1 1 2 7 5 1 0 2 4 4 2 25 1 8 1 -3 -1 5 2 12 2 23 -1 2 2 6 9 -7 5 6 6 1 5 9 5 2 2 1 4 1 2 -2 9 6 2 5 -1 5 -2 -2 2 2 6 13 0 2 8 5 2 -3 0 5 58 7 4 1 3 1 6 -1 9 0 5 10 12 2 9 2 10 8 9 8 5 2 1 4 3 4 8 6 1 3 14 3 7 3 -3 7 5 2 1 7 5 2 3 4 0 4 5 0 4 5 5 2 2 11 -1 13 18 6 9 9 6 2 7 3 12 11 0 32 2 6 5 4 2 3 1 6 5 4 1 4 4 5 3 4 1 14 4 12 2 5 2 11 4 6 0 4 3 0 8 0 1 2 5 10 3 1 6 8 0 2 2 6 4 4 6 13 3 0 5 3 2 5 2 7 4 4 2 1 -3 2 3 2 3 7 -3 4 5 8 1 3 7 1 1 -1 2 0 5 -2 3 1 9 0 1 -1 1 0 -1 0 3 1 9 0 2 2 2 2 2 1 6 5 7 4 2 -1 4 6 10 19 -1 7 14 -4 3 5 2 2 -1 0 1'

$ws = $wb.Worksheets.Item('YDS')
$ws.Range('C3').Value = '38 1 21 8 5 3 10 3 17 18 4 22 3 12 36 3 0 10 26 11 9 51 4 9 4 5 16 63 11 8 12 8 1 14 10 11 1 7 68 3 7 3 27 10 8 8 7 12 4 18 4 3 11 4 5 27 36 18 3 4 9 3 9 -4 13 30 8 6 3 6 15 6 54 5 29 3 53 4 0 9 29 11 8 5 18 5 -2 25 7 22 13 7 13 7 18 58 9 4 8 4 5 12 19 7 8 9 12 12 6 4 3 11 6 7 3 3 21 19 2 14 5 17 12 29 1 5 9 11 18 9 15 3 4 1 6 6 31 12 7 20 7 8 6 6 5 10 9 23 8 11 24 12 14 6 4 9 10 4 2 37 5 8 0 7 10 7 7 7 8 9 7 12 5 6 1 9 6 5 -2 11 3 9 38 8 8 14 12 4 0 7 1 10 6 18 20 26 15 6 17 6 3 4 9 22 13 7 2 13 3 16 5 9 6 9 3 10 8 6 7 6 46 11 3 12 21 21 8 10 34 26 22 5 17 15 8 10 9 8 12 3 13 10 5 6 2 8 7 6 11 4 22 20 16 13 12 38 5 19 5 3 21 10 7 20 9 11 5 41 4 7 7 8 14 3 7 5 7 3 12 3 31 7 14 6 5 5 9 9 31 41 6 6 3 2'

$ws = $wb.Worksheets.Item('ST')
$ws.Range('B4').Value = '61 61 55 62 52 68 66 63 63 59 57 65 68 67 65 65 65 66 59 64 61 64 60 69 67 50 57 65 64 66 66'

$ws = $wb.Worksheets.Item('ST')
$ws.Range('B5').Value = '25 7 18 16 17 24 26 16 15 24 16 24 21 30 18 18 27 31 15 19 21 13 27 23 26 17 17 21 26 37 16'

$ws = $wb.Worksheets.Item('ST')
$ws.Range('B6').Value = '14 19 14 19 13 19 16 15 16 12 17 17 20 17 21 31 17 17 25 25 27 20 7 14 26 17 22 20 22 30 16 25'

$ws = $wb.Worksheets.Item('ST')
$ws.Range('D3').Value = '47 59 45 41 53 42 12 63 30 42 41 39 45 44 42 63 44 45 49 50 52 37 49 41 57 35 42 42 36 41 42 46 40 27 53 47 48 35 37 35 39 46'

$ws = $wb.Worksheets.Item('ST')
$ws.Range('D4').Value = '11 13 0 0 0 0 0 14 0 18 0 0 0 0 0 6 14 0 0 11 16 0 0 0 0 0 0 7 0 4 0 0 0 0 10 15 8 0 0 0 0 0'

$ws = $wb.Worksheets.Item('ST')
$ws.Range('D5').Value = '0 0 0 0 0 0 10 13 9 12 18 0 14 16 22 0 0 0 6 0 0 2 0 0 0 0 2 19 0 0 11 0 0 0 0 9 0 0 0 16 0 4 0 0 0 2 0 0 1 10'

# --- Update numeric totals ---
$ws = $wb.Worksheets.Item('OFF')
$ws.Range('B2').Value = 6
$ws.Range('C2').Value = 196
$ws.Range('D2').Value = 10
$ws.Range('F2').Value = 67
$ws.Range('G2').Value = 46
$ws.Range('I2').Value = 4
$ws.Range('J2').Value = 27
$ws.Range('L2').Value = 228
$ws.Range('M2').Value = 152
$ws.Range('Q2').Value = 444
$ws.Range('C3').Value = 113
$ws.Range('E3').Value = 33
$ws.Range('F3').Value = 86
$ws.Range('G3').Value = 18
$ws.Range('H3').Value = 25
$ws.Range('I3').Value = 48
$ws.Range('J3').Value = 49
$ws.Range('N3').Value = 18
$ws = $wb.Worksheets.Item('DEF')
$ws.Range('C2').Value = 128
$ws.Range('F2').Value = 43
$ws.Range('G2').Value = 47
$ws.Range('I2').Value = 10
$ws.Range('L2').Value = 250
$ws.Range('M2').Value = 153
$ws.Range('O2').Value = 17
$ws.Range('P2').Value = 10
$ws.Range('Q2').Value = 422
$ws.Range('B3').Value = 15
$ws.Range('C3').Value = 149
$ws.Range('E3').Value = 33
$ws.Range('F3').Value = 81
$ws.Range('G3').Value = 24
$ws.Range('H3').Value = 17
$ws.Range('I3').Value = 49
$ws.Range('J3').Value = 48
$ws.Range('N3').Value = 16
$ws = $wb.Worksheets.Item('ST')
$ws.Range('B2').Value = 69
$ws.Range('D2').Value = 42
$ws.Range('F2').Value = 259
$ws.Range('G2').Value = 246
$ws.Range('J2').Value = 120
$ws.Range('K2').Value = 113
$ws.Range('L2').Value = 80
$ws.Range('M2').Value = 62
$ws.Range('B3').Value = 38
$ws = $wb.Worksheets.Item('TURNS')
$ws.Range('C2').Value = 7
$ws.Range('D3').Value = 7
$ws = $wb.Worksheets.Item('PEN')
$ws.Range('B3').Value = 9

$wb.Save()